$p = $ppt.ActivePresentation
Write-Output ("Path: " + $p.Path)
Write-Output ("FullName: " + $p.FullName)
Write-Output ("Name: " + $p.Name)
